# edit.ps1
# Applies the "Atualizacao de bases das ligas" commit:
#  1) A set of played-match rows (identified by same Date + clustered match ids)
#     had their records reshuffled/reordered within their date group. The row
#     number (col A), league (col C) and date (col D) stay put per physical
#     row, but everything else (id, teams, scores, odds...) moves to reflect
#     the new ordering.
#  2) A handful of still-unplayed fixtures (rows 326-333) got refreshed odds
#     (columns M/N/O/P/Q/R/T/U) - individual cell updates, no row reordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: reshuffle played-match rows within each date-cluster ---
# Read every affected row (columns B:AB) first so later writes never clobber
# a value that still needs to be read.
$row156 = $ws.Range("B156:AB156").Value2
$row157 = $ws.Range("B157:AB157").Value2
$row175 = $ws.Range("B175:AB175").Value2
$row176 = $ws.Range("B176:AB176").Value2
$row177 = $ws.Range("B177:AB177").Value2
$row180 = $ws.Range("B180:AB180").Value2
$row181 = $ws.Range("B181:AB181").Value2
$row182 = $ws.Range("B182:AB182").Value2
$row183 = $ws.Range("B183:AB183").Value2
$row184 = $ws.Range("B184:AB184").Value2
$row187 = $ws.Range("B187:AB187").Value2
$row188 = $ws.Range("B188:AB188").Value2
$row228 = $ws.Range("B228:AB228").Value2
$row229 = $ws.Range("B229:AB229").Value2
$row294 = $ws.Range("B294:AB294").Value2
$row295 = $ws.Range("B295:AB295").Value2
$row305 = $ws.Range("B305:AB305").Value2
$row306 = $ws.Range("B306:AB306").Value2
$row312 = $ws.Range("B312:AB312").Value2
$row313 = $ws.Range("B313:AB313").Value2

# Now write each row its new content per the permutation observed in the diff.
$ws.Range("B156:AB156").Value2 = $row157
$ws.Range("B157:AB157").Value2 = $row156
$ws.Range("B175:AB175").Value2 = $row177
$ws.Range("B176:AB176").Value2 = $row175
$ws.Range("B177:AB177").Value2 = $row176
$ws.Range("B180:AB180").Value2 = $row182
$ws.Range("B181:AB181").Value2 = $row180
$ws.Range("B182:AB182").Value2 = $row181
$ws.Range("B183:AB183").Value2 = $row188
$ws.Range("B184:AB184").Value2 = $row187
$ws.Range("B187:AB187").Value2 = $row184
$ws.Range("B188:AB188").Value2 = $row183
$ws.Range("B228:AB228").Value2 = $row229
$ws.Range("B229:AB229").Value2 = $row228
$ws.Range("B294:AB294").Value2 = $row295
$ws.Range("B295:AB295").Value2 = $row294
$ws.Range("B305:AB305").Value2 = $row306
$ws.Range("B306:AB306").Value2 = $row305
$ws.Range("B312:AB312").Value2 = $row313
$ws.Range("B313:AB313").Value2 = $row312

# --- Part 2: refresh odds on the upcoming-fixture rows (no reordering) ---
$ws.Range("M326").Value2 = 1.75
$ws.Range("O326").Value2 = 4.5
$ws.Range("Q326").Value2 = 2
$ws.Range("R326").Value2 = 1.85
$ws.Range("T326").Value2 = 1.85
$ws.Range("U326").Value2 = 2
$ws.Range("M327").Value2 = 1.2
$ws.Range("N327").Value2 = 5.75
$ws.Range("O327").Value2 = 11
$ws.Range("P327").Value2 = -2
$ws.Range("Q327").Value2 = 1.975
$ws.Range("R327").Value2 = 1.875
$ws.Range("T327").Value2 = 1.975
$ws.Range("U327").Value2 = 1.875
$ws.Range("M329").Value2 = 1.6
$ws.Range("N329").Value2 = 3.6
$ws.Range("O329").Value2 = 5.5
$ws.Range("Q329").Value2 = 1.8
$ws.Range("R329").Value2 = 2.05
$ws.Range("T329").Value2 = 1.875
$ws.Range("U329").Value2 = 1.975
$ws.Range("Q330").Value2 = 1.975
$ws.Range("R330").Value2 = 1.875
$ws.Range("M331").Value2 = 1.3
$ws.Range("O331").Value2 = 9
$ws.Range("Q331").Value2 = 2
$ws.Range("R331").Value2 = 1.85
$ws.Range("T331").Value2 = 2.025
$ws.Range("U331").Value2 = 1.825
$ws.Range("P332").Value2 = -0.75
$ws.Range("Q332").Value2 = 1.775
$ws.Range("R332").Value2 = 2.1
$ws.Range("M333").Value2 = 2.625
$ws.Range("O333").Value2 = 2.55
$ws.Range("P333").Value2 = 0
$ws.Range("Q333").Value2 = 2
$ws.Range("R333").Value2 = 1.85
